$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.370.38"
$ws.Range("E2").Value = "  +1.84%  "
$ws.Range("D3").Value = "1.826.62"
$ws.Range("E3").Value = "  +1.04%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'313.24"
$ws.Range("E5").Value = "  +0.98%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "'0.4468"
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("D8").Value = "'0.3767"
$ws.Range("E8").Value = "  +2.50%  "
$ws.Range("D9").Value = "'0.07415"
$ws.Range("E9").Value = "  +0.92%  "
$ws.Range("D10").Value = "'0.8784"
$ws.Range("E10").Value = "  +2.62%  "
$ws.Range("D11").Value = "'20.88"
$ws.Range("E11").Value = "  +0.99%  "
$ws.Range("D12").Value = "1.831.49"
$ws.Range("E12").Value = "  +1.30%  "
$ws.Range("D13").Value = "'6.713"
$ws.Range("E13").Value = "  +1.57%  "
$ws.Range("E14").Value = "  +2.22%  "
$ws.Range("D15").Value = "'93.04"
$ws.Range("D16").Value = "'0.07072"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").Value = "'0.000008827"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("D20").Value = "'15.10"
$ws.Range("E20").Value = "  +1.44%  "
$ws.Range("D21").Value = "27.373.22"
$ws.Range("E21").Value = "  +1.79%  "
$ws.Range("D22").Value = "'5.350"
$ws.Range("E22").Value = "  +3.69%  "
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("D24").Value = "'1.959"
$ws.Range("E24").Value = "  -1.69%  "
$ws.Range("D25").Value = "'151.22"
$ws.Range("E25").Value = "  -0.34%  "
$ws.Range("D26").Value = "'2.282"
$ws.Range("E26").Value = "  +4.04%  "
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("D28").Value = "'5.353"
$ws.Range("E28").Value = "  +2.51%  "
$ws.Range("D29").Value = "'117.28"
$ws.Range("E29").Value = "  +0.58%  "
$ws.Range("D30").Value = "'0.08914"
$ws.Range("E30").Value = "  +0.81%  "
$ws.Range("D31").Value = "'0.7951"
$ws.Range("D32").Value = "'1.198"
$ws.Range("E32").Value = "  +1.86%  "
$ws.Range("D33").Value = "'4.557"
$ws.Range("E33").Value = "  +1.86%  "
$ws.Range("D34").Value = "'2.962"
$ws.Range("E34").Value = "  +0.99%  "
$ws.Range("D35").Value = "'0.9998"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").Value = "'1.106"
$ws.Range("E36").Value = "  +1.45%  "
$ws.Range("E37").Value = "  +0.62%  "
$ws.Range("D38").Value = "'0.05279"
$ws.Range("E38").Value = "  +1.48%  "
$ws.Range("D39").Value = "'7.360"
$ws.Range("E39").Value = "  +4.76%  "
$ws.Range("D40").Value = "'0.5348"
$ws.Range("E40").Value = "  +0.34%  "
$ws.Range("D41").Value = "'2.876"
$ws.Range("E41").Value = "  +0.36%  "
$ws.Range("D42").Value = "'2.341"
$ws.Range("E42").Value = "  +18.55%  "
$ws.Range("D43").Value = "'0.1704"
$ws.Range("E43").Value = "  +0.72%  "
$ws.Range("D44").Value = "'8.695"
$ws.Range("E44").Value = "  +2.76%  "
$ws.Range("D45").Value = "'0.5095"
$ws.Range("E45").Value = "  -1.36%  "
$ws.Range("D46").Value = "'10.57"
$ws.Range("E46").Value = "  +0.52%  "
$ws.Range("E47").Value = "  -0.10%  "
$ws.Range("E48").Value = "  +0.62%  "
$ws.Range("D49").Value = "'0.9996"
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").Value = "'0.06389"
$ws.Range("E50").Value = "  +0.77%  "
$ws.Range("D51").Value = "'66.18"
$ws.Range("E51").Value = "  +5.66%  "
